# Restructure the "Input" sheet to the new standard template column layout,
# and drop the now-unused "I" column placeholder cells on the "갑지"/"을지" sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# ---- capture existing data before we clear the sheet ----
$numDataRows = 7   # rows 2..8

$company      = $ws.Cells.Item(2, 1).Value2   # A (거래처명) - same for every row
$site         = $ws.Cells.Item(2, 2).Value2   # B (현장명) - same for every row
$companyEmail = "티에스이앤씨@example.com"
$deliveryEmail = "delivery@example.com"

$poDate  = @()
$dueDate = @()
$item    = @()
$spec    = @()
$qty     = @()
$price   = @()
$total   = @()
$cat1    = @()
$cat2    = @()

for ($i = 0; $i -lt $numDataRows; $i++) {
    $r = 2 + $i
    $poDate  += $ws.Cells.Item($r, 3).Value2    # C 발주일
    $dueDate += $ws.Cells.Item($r, 4).Value2    # D 납기일
    $item    += $ws.Cells.Item($r, 6).Value2    # F 품목
    $spec    += $ws.Cells.Item($r, 7).Value2    # G 규격
    $qty     += $ws.Cells.Item($r, 8).Value2    # H 수량
    $price   += $ws.Cells.Item($r, 10).Value2   # J 단가
    $total   += $ws.Cells.Item($r, 13).Value2   # M 합계
    $cat1    += $ws.Cells.Item($r, 14).Value2   # N 대분류
    $cat2    += $ws.Cells.Item($r, 15).Value2   # O 중분류
}

# ---- clear the whole sheet (data + header styling) ----
$ws.Cells.Clear()

# ---- write the new header row (no bold/border style) ----
$headers = @("발주일자","납기일자","거래처명","거래처 이메일","납품처명","납품처 이메일","프로젝트명","대분류","중분류","소분류","품목명","규격","수량","단가","총금액","비고")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value2 = $headers[$c]
}

# ---- write the new data rows in the new column order ----
for ($i = 0; $i -lt $numDataRows; $i++) {
    $r = 2 + $i

    # The 발주일자/납기일자 values look like dates ("2025-09-10"), and Excel
    # auto-converts such text to a date serial on assignment. Temporarily
    # mark the cell as Text first so the original literal string is kept,
    # then restore the Normal style so no stray formatting is introduced.
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value2 = $poDate[$i]          # A 발주일자
    $ws.Cells.Item($r, 1).Style = "Normal"

    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value2 = $dueDate[$i]         # B 납기일자
    $ws.Cells.Item($r, 2).Style = "Normal"

    $ws.Cells.Item($r, 3).Value2  = $company            # C 거래처명
    $ws.Cells.Item($r, 4).Value2  = $companyEmail        # D 거래처 이메일
    $ws.Cells.Item($r, 5).Value2  = $site                # E 납품처명
    $ws.Cells.Item($r, 6).Value2  = $deliveryEmail        # F 납품처 이메일
    $ws.Cells.Item($r, 7).Value2  = $site                # G 프로젝트명
    $ws.Cells.Item($r, 8).Value2  = $cat1[$i]            # H 대분류
    $ws.Cells.Item($r, 9).Value2  = $cat2[$i]            # I 중분류
    # column J (소분류) intentionally left blank
    $ws.Cells.Item($r, 11).Value2 = $item[$i]            # K 품목명
    $ws.Cells.Item($r, 12).Value2 = $spec[$i]            # L 규격
    $ws.Cells.Item($r, 13).Value2 = $qty[$i]             # M 수량
    $ws.Cells.Item($r, 14).Value2 = $price[$i]           # N 단가
    $ws.Cells.Item($r, 15).Value2 = $total[$i]           # O 총금액
    # column P (비고) intentionally left blank
}

# ---- 갑지 / 을지: drop the leftover empty "I" placeholder cells ----
foreach ($sheetName in @("갑지", "을지")) {
    $sh = $wb.Worksheets.Item($sheetName)
    for ($r = 2; $r -le 8; $r++) {
        $sh.Cells.Item($r, 9).ClearContents()
    }
}
